# Indicateurs2.xlsx - "Add files via upload" edit
# Replaces the 2-column PA_U1/PA_U2 mini-table with an 8-column
# PO_*/AOSO_* table (U1/U2/U3/LPV indicators) and restyles the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Propagate existing formats to their new homes BEFORE the source
#    cells' own values/styles get overwritten later in the script.
# ---------------------------------------------------------------------

# A2's current format (the shaded/centered "value" style) is needed on
# C2 and E2 in the new layout - grab it first.
$ws.Range("A2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)

# A1's current format (plain horizontal-center style) is needed on the
# rest of the numeric row (A2, B2, D2, F2, G2, H2).
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Header row: new label text, written in the exact order that
#    reproduces the target shared-strings table, then a fresh
#    centered/middle style built on G1 and copied across the row.
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "PO_U1"
$ws.Range("E1").Value = "PO_U2"
$ws.Range("G1").Value = "PO_U3"
$ws.Range("A1").Value = "PO_LPV"
$ws.Range("D1").Value = "AOSO_U1"
$ws.Range("F1").Value = "AOSO_U2"
$ws.Range("H1").Value = "AOSO_U3"
$ws.Range("B1").Value = "AOSO_LPV"

$ws.Range("G1").HorizontalAlignment = -4108
$ws.Range("G1").VerticalAlignment = -4108
$ws.Range("G1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Data row values.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 33
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 75
$ws.Range("E2").Value = 26
$ws.Range("F2").Value = 46
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 32

# ---------------------------------------------------------------------
# 4) Selection, matching the saved workbook's cursor position.
# ---------------------------------------------------------------------
$ws.Range("F9").Select() | Out-Null

Write-Output "Indicateurs2 table rebuilt (A1:H2)"
